# Generate Report for Handoff
#
# The b4a3cc72-6256-4344-84a5-c8a9f73aa3ee file has finished its handback
# cycle and is ready to be sent out again, so its status flips from
# "Handed back: in sync with en-US" to "Ready for handoff" everywhere it is
# reported, and the associated handoff timestamps are bumped accordingly.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet --------------------------------------------------
# Columns: A=File Name, B=zh-cn, C=de-de, D=Latest Handoff Date
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$overview.Range("D3").Value = "2016-39-21 02:39:25"

# --- zh-cn sheet -------------------------------------------------------
# Columns: C=Status, E=Latest Handoff Datetime
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("E3").Value = "2016-03-21 02:39:21"

# --- de-de sheet -------------------------------------------------------
# Columns: C=Status, E=Latest Handoff Datetime
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("E3").Value = "2016-03-21 02:39:25"
